$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.45, 0.266, -0.508, 36, 47, 139),
    @(-0.06, 0.404, -0.304, -178, 12, 165),
    @(-0.076, -0.138, 0.173, 100, 74, -120),
    @(0.033, -0.114, -0.482, -51, 31, 153),
    @(0.164, 0.078, 0.163, -179, 22, 14),
    @(-0.116, -0.072, -0.288, -12, 0, 145),
    @(-0.343, -0.052, 0.511, -175, 77, -28),
    @(-0.079, 0.334, 0.028, 42, 31, 125),
    @(-0.121, -0.307, -0.285, -29, -60, -134),
    @(-0.384, -0.084, 0.701, 96, 12, -105)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($i + 1, $j + 1).Value = $rowVals[$j]
    }
}
